$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo in student email: "pasne.d@husky.neu.edu" -> "panse.d@husky.neu.edu"
$ws.Range("C2").Value = "panse.d@husky.neu.edu"

# Update the active selection to reflect where the user clicked after the edit
$ws.Range("C2").Select()
